$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Question1" entry was inserted as the first School-category question,
# which pushes every subsequent question/category/answer row down by one.
# Rewrite rows 2-16 (A:I) in place to reflect the shifted data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "School"
$ws.Range("D2").Value = "Question1"
$ws.Range("E2").Value = "Answer1"
$ws.Range("F2").Value = "Answer2"
$ws.Range("G2").Value = "Answer3"
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = "Answer2"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "School"
$ws.Range("D3").Value = "Question2"
$ws.Range("E3").Value = "Answer1"
$ws.Range("F3").Value = "Answer2"
$ws.Range("G3").Value = "Answer3"
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = "Answer1"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "School"
$ws.Range("D4").Value = "Question3"
$ws.Range("E4").Value = "Answer1"
$ws.Range("F4").Value = "Answer2"
$ws.Range("G4").Value = "Answer3"
$ws.Range("H4").Value = 7
$ws.Range("I4").Value = "Answer3"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "School"
$ws.Range("D5").Value = "Question4"
$ws.Range("E5").Value = "Answer1"
$ws.Range("F5").Value = "Answer2"
$ws.Range("G5").Value = "Answer3"
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = "Answer1"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "School"
$ws.Range("D6").Value = "Question5"
$ws.Range("E6").Value = "Answer1"
$ws.Range("F6").Value = "Answer2"
$ws.Range("G6").Value = "Answer3"
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = "Answer2"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = "Home"
$ws.Range("D7").Value = "Question6"
$ws.Range("E7").Value = "Answer1"
$ws.Range("F7").Value = "Answer2"
$ws.Range("G7").Value = "Answer3"
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = "Answer2"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = "Home"
$ws.Range("D8").Value = "Question7"
$ws.Range("E8").Value = "Answer1"
$ws.Range("F8").Value = "Answer2"
$ws.Range("G8").Value = "Answer3"
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Answer2"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = "Home"
$ws.Range("D9").Value = "Question8"
$ws.Range("E9").Value = "Answer1"
$ws.Range("F9").Value = "Answer2"
$ws.Range("G9").Value = "Answer3"
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = "Answer3"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = "Home"
$ws.Range("D10").Value = "Question9"
$ws.Range("E10").Value = "Answer1"
$ws.Range("F10").Value = "Answer2"
$ws.Range("G10").Value = "Answer3"
$ws.Range("H10").Value = 14
$ws.Range("I10").Value = "Answer2"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = "Home"
$ws.Range("D11").Value = "Question10"
$ws.Range("E11").Value = "Answer1"
$ws.Range("F11").Value = "Answer2"
$ws.Range("G11").Value = "Answer3"
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = "Answer1"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = "Public Places"
$ws.Range("D12").Value = "Question11"
$ws.Range("E12").Value = "Answer1"
$ws.Range("F12").Value = "Answer2"
$ws.Range("G12").Value = "Answer3"
$ws.Range("H12").Value = 8
$ws.Range("I12").Value = "Answer2"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = "Public Places"
$ws.Range("D13").Value = "Question12"
$ws.Range("E13").Value = "Answer1"
$ws.Range("F13").Value = "Answer2"
$ws.Range("G13").Value = "Answer3"
$ws.Range("H13").Value = 12
$ws.Range("I13").Value = "Answer3"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = "Public Places"
$ws.Range("D14").Value = "Question13"
$ws.Range("E14").Value = "Answer1"
$ws.Range("F14").Value = "Answer2"
$ws.Range("G14").Value = "Answer3"
$ws.Range("H14").Value = 9
$ws.Range("I14").Value = "Answer1"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = "Public Places"
$ws.Range("D15").Value = "Question14"
$ws.Range("E15").Value = "Answer1"
$ws.Range("F15").Value = "Answer2"
$ws.Range("G15").Value = "Answer3"
$ws.Range("H15").Value = 7
$ws.Range("I15").Value = "Answer1"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = "Public Places"
$ws.Range("D16").Value = "Question15"
$ws.Range("E16").Value = "Answer1"
$ws.Range("F16").Value = "Answer2"
$ws.Range("G16").Value = "Answer3"
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = "Answer2"

# Row 16 is brand new (previously the sheet stopped at row 15), so give its
# index cell (column A) the same bold/centered/bordered formatting used by
# every other cell in column A, and match the row height used elsewhere.
$a16 = $ws.Range("A16")
$a16.Font.Bold = $true
$a16.HorizontalAlignment = -4108
$a16.VerticalAlignment = -4160
$a16.Borders.Item(7).LineStyle = 1
$a16.Borders.Item(8).LineStyle = 1
$a16.Borders.Item(9).LineStyle = 1
$a16.Borders.Item(10).LineStyle = 1
$ws.Rows.Item(16).RowHeight = 15

# New trailing blank row 17 with only a formatted (but empty) column-A cell,
# matching the look of the index column used throughout the table.
$a17 = $ws.Range("A17")
$a17.Font.Bold = $true
$a17.HorizontalAlignment = -4108
$a17.VerticalAlignment = -4160
$a17.Borders.Item(7).LineStyle = 1
$a17.Borders.Item(8).LineStyle = 1
$a17.Borders.Item(9).LineStyle = 1
$a17.Borders.Item(10).LineStyle = 1
$ws.Rows.Item(17).RowHeight = 15

# Resize columns to fit the new content (mirrors the bestFit column widths Excel computed)
$ws.Columns.Item(2).ColumnWidth = 5.035714285714286
$ws.Columns.Item(3).ColumnWidth = 11.285714285714286
$ws.Columns.Item(4).ColumnWidth = 9.535714285714286
$ws.Columns.Item(5).ColumnWidth = 8.785714285714285
$ws.Columns.Item(6).ColumnWidth = 9.071428571428571
$ws.Columns.Item(7).ColumnWidth = 9.071428571428571
$ws.Columns.Item(9).ColumnWidth = 11.928571428571427

# Update the active selection to the newly edited cell
$ws.Range("D2").Select()
